$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AlgoHive")

# --- Scoring inputs (top table) ---
$ws.Range("H2").Value = "MEDIUM"
$ws.Range("H4").Value = 9
$ws.Range("H5").Value = "<20"
$ws.Range("H6").Value = 0

# --- Totals section: B50 becomes a formula, B51 range extended ---
$ws.Range("B50").Formula = "=SUM(B34:B40)"
$ws.Range("B51").Formula = "=SUM(B34:B44)"

# --- Bottom grade table ---
$ws.Range("B55").Value = 800
$ws.Range("C55").Formula = "=MIN(20, ROUNDUP(IF(B55<=B50, (B55/B50)*10, 10 + ((B55-B50)/(B51-B50))*10), 0.5))"

# --- View state ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("F58").Select()
